# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P columns (Late / heading / Outstanding) one
#   column to the right (-> O/P/Q).
# - Make "Repayment schedule" the active sheet/tab, with selection on P15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at position 14 (column N); everything from N onward
# (including styles/values) shifts right automatically.
$ws.Columns.Item(14).Insert()

# New column width for the inserted column N (stored width "11").
$ws.Columns.Item(14).ColumnWidth = 10.2

# Make "Repayment schedule" the active sheet/tab and set its selection.
$ws.Activate()
$ws.Range("P15").Select()
